$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (G1) onto the new
# "Save" header cell (H1), then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data column values (era data update).
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
